$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ D = "26.887.97"; E = "  -0.94%  " }
    3  = @{ D = "1.615.27"; E = "  -1.34%  " }
    5  = @{ D = "212.11"; E = "  -2.29%  " }
    6  = @{ E = "  -1.11%  " }
    7  = @{ E = "  -0.03%  " }
    8  = @{ E = "  -1.82%  " }
    10 = @{ D = "19.77"; E = "  -1.69%  " }
    11 = @{ D = "0.0837"; E = "  -1.32%  " }
    12 = @{ D = "1.843.99"; E = "  -1.26%  " }
    13 = @{ D = "1.610.64"; E = "  -1.72%  " }
    14 = @{ E = "  -1.33%  " }
    15 = @{ D = "0.533"; E = "  -2.07%  " }
    16 = @{ D = "26.898.09"; E = "  -0.96%  " }
    17 = @{ D = "63.80" }
    18 = @{ D = "0.0₃0730"; E = "  -0.93%  " }
    19 = @{ D = "210.78"; E = "  -2.68%  " }
    20 = @{ E = "  +0.04%  " }
    21 = @{ D = "6.75"; E = "  -1.84%  " }
    22 = @{ D = "4.30"; E = "  -2.60%  " }
    23 = @{ E = "  -7.35%  " }
    24 = @{ D = "8.92"; E = "  -2.34%  " }
    25 = @{ D = "146.38"; E = "  -0.86%  " }
    26 = @{ E = "  +1.09%  " }
    27 = @{ E = "  +0.00%  " }
    28 = @{ D = "0.113"; E = "  -4.32%  " }
    29 = @{ D = "15.39"; E = "  -1.70%  " }
    30 = @{ D = "0.0503"; E = "  -1.07%  " }
    31 = @{ E = "  -1.92%  " }
    32 = @{ D = "3.26"; E = "  -3.37%  " }
    33 = @{ D = "0.701"; E = "  +27.69%  " }
    34 = @{ D = "2.96"; E = "  -2.13%  " }
    35 = @{ D = "1.321.33"; E = "  +1.27%  " }
    36 = @{ E = "  -2.25%  " }
    37 = @{ E = "  -0.56%  " }
    38 = @{ E = "  -1.53%  " }
    39 = @{ D = "0.828"; E = "  -2.90%  " }
    40 = @{ E = "  -0.02%  " }
    41 = @{ E = "  -2.00%  " }
    42 = @{ E = "  -2.61%  " }
    43 = @{ E = "  -1.19%  " }
    44 = @{ D = "63.42"; E = "  +1.90%  " }
    45 = @{ D = "1.751.32"; E = "  -1.52%  " }
    46 = @{ D = "89.38"; E = "  -1.64%  " }
    47 = @{ D = "1.61"; E = "  +0.52%  " }
    48 = @{ D = "0.811"; E = "  +7.58%  " }
    49 = @{ D = "0.0₆0104"; E = "  -1.38%  " }
    50 = @{ D = "0.0514"; E = "  -0.25%  " }
    51 = @{ D = "0.0981"; E = "  +2.39%  " }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $cell = $ws.Range("$col$row")
        $cell.NumberFormat = "@"
        $cell.Value = $cols[$col]
    }
}
